# "Generate Report for Handoff"
#
# The d46f061f-73de-4abe-93a9-1c8b5d4dca03 file moved from
# "Handed back: in sync with en-US" status to "Ready for handoff",
# with refreshed handoff timestamps, on the Overview sheet and on
# both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the d46f061f-...-93a9-1c8b5d4dca03.md file ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-33-14 08:33:19"

# --- zh-cn sheet: row 3 is the same file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-14 08:33:16"

# --- de-de sheet: row 3 is the same file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-14 08:33:19"
